$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-35 (Generation 0-33) -> Fitness 7594
$ws.Range("C2:C35").Value = 7594

# Rows 36-252 (Generation 34-250) -> Fitness 7586
$ws.Range("C36:C252").Value = 7586
